# ResourceFile_Method_HT.xlsx - "Incidence of HT is (nearly) correct"
#
# The divisor used to convert prevalence -> incidence in column C of the
# "incidence2018_plus" sheet changes for three blocks of rows:
#   rows 37-46  : /6  -> /12
#   rows 47-56  : /6  -> /11
#   rows 57-122 : /6  -> /12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("incidence2018_plus")

for ($r = 37; $r -le 46; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/12"
}

for ($r = 47; $r -le 56; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/11"
}

for ($r = 57; $r -le 122; $r++) {
    $ws.Range("C$r").Formula = "=prevalence2018!C$r/12"
}

# Reflect the author's updated cursor/selection position on this sheet.
$ws.Activate()
$ws.Range("G41").Select()
